$wb = $excel.ActiveWorkbook

# The "Personnel" sheet holds the transect personnel info and needs the
# ORCID (userId) added for Diana Fontaine in row 11, column F.
$ws = $wb.Worksheets.Item("Personnel")

$ws.Range("F11").Value = "0000-0001-9172-6904"

# Move the active selection to mirror where editing ended up.
$ws.Range("F22").Select()
